$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.660.99"
$ws.Range("E2").Value = "  -1.51%  "

$ws.Range("D3").Value = "1.615.52"
$ws.Range("E3").Value = "  -1.88%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5065"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.93%  "

$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2556"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.80%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06340"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07768"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.228"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.23%  "

$ws.Range("D13").Value = "1.624.01"
$ws.Range("E13").Value = "  -1.54%  "

$ws.Range("D14").Value = "1.838.29"
$ws.Range("E14").Value = "  -1.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5536"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.62%  "

$ws.Range("D17").Value = "0.0₅7482"
$ws.Range("E17").Value = "  -3.95%  "

$ws.Range("D18").Value = "25.669.88"
$ws.Range("E18").Value = "  -1.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.333"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.718"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.943"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.22%  "

$ws.Range("E25").Value = "  -3.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.55%  "

$ws.Range("E27").Value = "  +2.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.700"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.231"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04837"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.282"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.166"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.542"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.365"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8881"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.81%  "

$ws.Range("D37").Value = "1.122.33"
$ws.Range("E37").Value = "  +0.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.526"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5460"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.80%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01552"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.72%  "

$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.002"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.552"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7908"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.26%  "

$ws.Range("D45").Value = "1.763.93"
$ws.Range("E45").Value = "  -1.15%  "

$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  -8.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4407"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05087"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.537"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.70%  "
